# Insert a new data row at row 212 (pushing existing rows 212..304 down to 213..305)
# and populate the new row with the latest weekly price observation for Perejil.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(212).Insert()

$ws.Range("A212").Value = 4
$ws.Range("B212").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C212").Value = "Los Lagos"
$ws.Range("D212").Value = 44875
$ws.Range("E212").Value = 10
$ws.Range("F212").Value = 100112044
$ws.Range("G212").Value = "Perejil"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 80
$ws.Range("K212").Value = 6000
$ws.Range("L212").Value = 6000
$ws.Range("M212").Value = 6000
$ws.Range("N212").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O212").Value = "Región de La Araucanía"
$ws.Range("P212").Value = 2000
$ws.Range("Q212").Value = 3
$ws.Range("R212").Value = "Hortaliza"
